# Update Sheet2 ("measurement level" style data):
#   - Header row stays textual (header labels changed to type names)
#   - Data rows: first and third columns become numeric (Double / integer),
#     middle column stays text (string)
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("B1").Value = "Double"
$ws2.Range("C1").Value = "string"
$ws2.Range("D1").Value = "integer"

$ws2.Range("B2").Value = 0.76
$ws2.Range("C2").Value = "str8"
$ws2.Range("D2").Value = 1700

$ws2.Range("B3").Value = 1.2
$ws2.Range("C3").Value = "str11"
$ws2.Range("D3").Value = 12

# Selection / active-tab bookkeeping: Sheet2 becomes the selected/active sheet
$ws2.Range("C22").Select()
$ws2.Activate()
